# "changed color to neural[sic: neutral]"
# The three accent rectangles that sit behind the neural-network diagram on
# slides 2 and 3 were recolored to the theme's neutral accent (Accent 5),
# replacing the hard-coded red (FF8F8F) / green (77E982) fills and
# stripping the extra luminance-mod/-off tint that had been applied to the
# rectangle that was already using Accent 5 -- all three shapes now share
# a plain Accent 5 fill at 50% transparency.

$p = $ppt.ActivePresentation

for ($slideIndex = 2; $slideIndex -le 3; $slideIndex++) {
    $slide = $p.Slides.Item($slideIndex)

    foreach ($shapeName in @("Rectangle 113", "Rectangle 155", "Rectangle 192")) {
        $shape = $slide.Shapes.Item($shapeName)

        # ppAccent5 == 9 in the PpThemeColorIndex enumeration.
        $shape.Fill.ForeColor.ObjectThemeColor = 9
        # Preserve the 50% transparency (alpha=50000) the fills already had.
        $shape.Fill.Transparency = 0.5
    }
}
